# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
# Updates the "Periodo Mora" detail table (rows 16-67) on sheet Hoja1 with the
# refreshed worker/period/value data for NIT 8904000631.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$data = @(
    @(16, "73545140", "EZEQUIEL ARROYO MARQUEZ", "1908", 24640, 781242),
    @(17, "73545140", "EZEQUIEL ARROYO MARQUEZ", "1907", 24640, 781242),
    @(18, "73545140", "EZEQUIEL ARROYO MARQUEZ", "1906", 24640, 781242),
    @(19, "73545140", "EZEQUIEL ARROYO MARQUEZ", "1905", 24640, 781242),
    @(20, "73545140", "EZEQUIEL ARROYO MARQUEZ", "1904", 24640, 781242),
    @(21, "73545140", "EZEQUIEL ARROYO MARQUEZ", "1903", 24640, 781242),
    @(22, "73545140", "EZEQUIEL ARROYO MARQUEZ", "1902", 24640, 781242),
    @(23, "73545140", "EZEQUIEL ARROYO MARQUEZ", "1901", 24640, 781242),
    @(24, "73545140", "EZEQUIEL ARROYO MARQUEZ", "1812", 24640, 781242),
    @(25, "73545140", "EZEQUIEL ARROYO MARQUEZ", "1811", 24640, 781242),
    @(26, "73545140", "EZEQUIEL ARROYO MARQUEZ", "1810", 24640, 781242),
    @(27, "73545140", "EZEQUIEL ARROYO MARQUEZ", "1809", 24640, 781242),
    @(28, "73545140", "EZEQUIEL ARROYO MARQUEZ", "1808", 24640, 781242),
    @(29, "73545140", "EZEQUIEL ARROYO MARQUEZ", "1807", 24640, 781242),
    @(30, "73545140", "EZEQUIEL ARROYO MARQUEZ", "1806", 24640, 781242),
    @(31, "73545140", "EZEQUIEL ARROYO MARQUEZ", "1805", 24640, 781242),
    @(32, "73545140", "EZEQUIEL ARROYO MARQUEZ", "1804", 24640, 781242),
    @(33, "73545140", "EZEQUIEL ARROYO MARQUEZ", "1803", 24640, 781242),
    @(34, "73545140", "EZEQUIEL ARROYO MARQUEZ", "1802", 24640, 781242),
    @(35, "73545140", "EZEQUIEL ARROYO MARQUEZ", "1801", 24640, 781242),
    @(36, "73545140", "EZEQUIEL ARROYO MARQUEZ", "1712", 24640, 781242),
    @(37, "73545140", "EZEQUIEL ARROYO MARQUEZ", "1711", 24640, 781242),
    @(38, "73545140", "EZEQUIEL ARROYO MARQUEZ", "1710", 24640, 781242),
    @(39, "73545140", "EZEQUIEL ARROYO MARQUEZ", "1709", 24640, 781242),
    @(40, "73545140", "EZEQUIEL ARROYO MARQUEZ", "1708", 24640, 781242),
    @(41, "73545140", "EZEQUIEL ARROYO MARQUEZ", "1707", 24640, 781242),
    @(42, "73545140", "EZEQUIEL ARROYO MARQUEZ", "1706", 31249, 781242),
    @(43, "73545140", "EZEQUIEL ARROYO MARQUEZ", "1705", 31249, 781242),
    @(44, "73545140", "EZEQUIEL ARROYO MARQUEZ", "1704", 31249, 781242),
    @(45, "73545140", "EZEQUIEL ARROYO MARQUEZ", "1703", 31249, 781242),
    @(46, "73545140", "EZEQUIEL ARROYO MARQUEZ", "1702", 31249, 781242),
    @(47, "73545140", "EZEQUIEL ARROYO MARQUEZ", "1701", 31249, 781242),
    @(48, "9173408", "MARCIAL RAFAEL DIAZ ALMEIDA", "1902", 60000, 1500000),
    @(49, "73545110", "ISMAEL ANTONIO HERNANDEZ ACOSTA", "1902", 60000, 1500000),
    @(50, "73545140", "EZEQUIEL ARROYO MARQUEZ", "1903", 31249, 781242),
    @(51, "9173408", "MARCIAL RAFAEL DIAZ ALMEIDA", "1903", 60000, 1500000),
    @(52, "73545110", "ISMAEL ANTONIO HERNANDEZ ACOSTA", "1903", 60000, 1500000),
    @(53, "73545140", "EZEQUIEL ARROYO MARQUEZ", "1904", 31249, 781242),
    @(54, "9173408", "MARCIAL RAFAEL DIAZ ALMEIDA", "1904", 60000, 1500000),
    @(55, "73545110", "ISMAEL ANTONIO HERNANDEZ ACOSTA", "1904", 60000, 1500000),
    @(56, "73545140", "EZEQUIEL ARROYO MARQUEZ", "1905", 31249, 781242),
    @(57, "9173408", "MARCIAL RAFAEL DIAZ ALMEIDA", "1905", 60000, 1500000),
    @(58, "73545110", "ISMAEL ANTONIO HERNANDEZ ACOSTA", "1905", 60000, 1500000),
    @(59, "73545140", "EZEQUIEL ARROYO MARQUEZ", "1906", 31249, 781242),
    @(60, "9173408", "MARCIAL RAFAEL DIAZ ALMEIDA", "1906", 60000, 1500000),
    @(61, "73545110", "ISMAEL ANTONIO HERNANDEZ ACOSTA", "1906", 60000, 1500000),
    @(62, "73545140", "EZEQUIEL ARROYO MARQUEZ", "1907", 31249, 781242),
    @(63, "9173408", "MARCIAL RAFAEL DIAZ ALMEIDA", "1907", 60000, 1500000),
    @(64, "73545110", "ISMAEL ANTONIO HERNANDEZ ACOSTA", "1907", 60000, 1500000),
    @(65, "73545140", "EZEQUIEL ARROYO MARQUEZ", "1908", 19791, 781242),
    @(66, "9173408", "MARCIAL RAFAEL DIAZ ALMEIDA", "1908", 38000, 1500000),
    @(67, "73545110", "ISMAEL ANTONIO HERNANDEZ ACOSTA", "1908", 38000, 1500000)
)

foreach ($item in $data) {
    $r = $item[0]
    $ws.Cells.Item($r, 3).Value = $item[1]   # C: N Doc Trabajador
    $ws.Cells.Item($r, 4).Value = $item[2]   # D: Nombre Trabajador
    $ws.Cells.Item($r, 5).Value = $item[3]   # E: Periodo Mora
    $ws.Cells.Item($r, 6).Value = $item[4]   # F: Valor Mora
    $ws.Cells.Item($r, 7).Value = $item[5]   # G: Salario Basico
}